$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'basketball leg sleeve youth padded'
    2 = 'knee pad gym'
    3 = 'work need pads'
    4 = 'softball material'
    5 = 'compression spandex men'
    6 = 'football leggings men'
    7 = 'knees bees'
    8 = 'snowboarding pants youth'
    9 = 'knee length pants'
    10 = 'medias de basketball'
    11 = 'padded sliding shorts youth'
    12 = 'men gym tight pants'
    13 = 'compression tight leggings'
    14 = 'hockey leggings girls'
    15 = 'boys black tight pants'
    16 = 'yoga pad'
    17 = 'leg compression pants men'
    18 = 'wrestling shorts men'
    19 = 'yoga pants men big and tall'
    20 = 'wrestling knee'
    21 = 'youth sliding shorts'
    22 = 'youth hockey pants'
    23 = 'sliding shorts baseball'
    24 = 'knee pads replacement'
    25 = 'thread protector 1/2 x 28'
    26 = 'running tights youth'
    27 = 'mesh leggings men'
    28 = 'boys hockey pants'
    29 = 'baseball pants youth boys black'
    30 = 'basketball tight shorts for boys'
    31 = 'sports compression leggings'
    32 = 'knee pads for yoga'
    33 = 'padded work pants mens'
    34 = 'baseball pants for boys'
    35 = 'youth spandex'
    36 = 'women knee pads for work'
    37 = 'big and tall mens compression pants'
    38 = 'padded compression shorts youth'
    39 = 'boys tights youth'
    40 = 'men sheer pants'
    41 = 'football girdle with pads for men'
    42 = 'knee pads for work men'
    43 = 'yoga position chart'
    44 = 'youth football leggings boys'
    45 = 'knee compression sleeve pad'
    46 = 'black baseball pants mens'
    47 = 'mens leggings shorts'
    48 = 'youth baseball pants knee high'
    49 = 'elastic waist baseball pants'
    50 = 'adult tights'
    51 = 'knee pads for work black'
    52 = 'arthritis test'
    53 = 'cycling sweat guard'
    54 = 'volleyball knee pads extra large'
    55 = 'size 5 basketball'
    56 = 'athletic knee compression'
    57 = 'knee pads for man'
    58 = 'youth basketball knee sleeve'
    59 = 'knee pads for mountain biking'
    60 = 'best knee pads for work'
    61 = 'hip pads for men'
    62 = 'girl compression pants'
    63 = 'compression shorts basketball'
    64 = 'basketball knee sleeve boys'
    65 = 'mens pants big and tall'
    66 = 'below knee shorts men'
    67 = 'youth wrestling shorts'
    68 = 'capri spandex'
    69 = 'yoga pads for hands'
    70 = 'football pants adult with pads'
    71 = 'boys youth compression pants'
    72 = 'exercise kneeling pad'
    73 = 'knee pad volleyball'
    74 = 'knee pads working'
    75 = 'baseball softball pants'
    76 = 'boys knee pads volleyball'
    77 = 'sliding workout pads'
    78 = 'knee pads for youth'
    79 = 'black legging for men'
    80 = 'cheap leggings for men'
    81 = 'little boys compression leggings'
    82 = 'basketball compression gear'
    83 = 'men sport pants'
    84 = 'cold knee pad'
    85 = 'black compression shorts for men'
    86 = 'calf sleeves for men basketball'
    87 = 'mens work knee pads'
    88 = 'the bees knees'
    89 = 'athletic leggings for men'
    90 = 'compression pants for boys'
    91 = 'patella knee pads'
    92 = 'lacrosse tights'
    93 = 'boys leggings youth'
    94 = 'boy sport tights'
    95 = 'volleyball pants'
    96 = 'knee pads for wrestling'
    97 = 'football padded shorts for men'
    98 = 'mens basketball outdoor'
    99 = 'compression pants size'
    100 = 'mens knee pads construction'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
